# feat: add 2022-Q4 data
#
# The workbook currently has 3 sheets: 总计 (totals), 2022-Q1, 2021-Q4.
# This script inserts a brand new "2022-Q4" sheet (with its own fund-holding
# table) positioned right after "总计" and before "2022-Q1", and updates the
# "总计" summary sheet with a new row describing the 2022-Q4 snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Re-fetch the "2022-Q1" sheet by name (worksheet indices shifted because of
# the insert above) and copy its header row + first data row into the new
# sheet so the new sheet starts out with the same column headers/styles.
# (Column A of row 1 is intentionally left untouched - the source template
# has no A1 cell either.)
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Range("B1:H1").Copy($q4.Range("B1:H1"))
$q1.Range("A2:H2").Copy($q4.Range("A2:H2"))

# Replicate the formatted data row down to rows 3-6 (5 fund rows total).
$q4.Range("A2:H2").Copy($q4.Range("A3:H3"))
$q4.Range("A2:H2").Copy($q4.Range("A4:H4"))
$q4.Range("A2:H2").Copy($q4.Range("A5:H5"))
$q4.Range("A2:H2").Copy($q4.Range("A6:H6"))

# Columns B (fund code) and D/E/F/G (scale/position/weight/value) are stored
# as plain text in the source data (fund codes have leading zeros, and the
# decimal values keep a fixed number of trailing zeros) - force text format
# so the literal strings survive instead of being coerced to numbers.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

# ---------------------------------------------------------------------
# 2) Fill in the 2022-Q4 fund-holding data.
# ---------------------------------------------------------------------

# Row 2 - 007484 信澳核心科技混合
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "007484"
$q4.Range("C2").Value = "信澳核心科技混合"
$q4.Range("D2").Value = "21.50"
$q4.Range("E2").Value = "93.51"
$q4.Range("F2").Value = "3.60"
$q4.Range("G2").Value = "0.7740"
$q4.Range("H2").Value = 5

# Row 3 - 003956 南方产业智选股票
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "003956"
$q4.Range("C3").Value = "南方产业智选股票"
$q4.Range("D3").Value = "3.60"
$q4.Range("E3").Value = "85.80"
$q4.Range("F3").Value = "4.78"
$q4.Range("G3").Value = "0.1721"
$q4.Range("H3").Value = 6

# Row 4 - 011214 招商惠润一年定期开放混合（MOM）A
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "011214"
$q4.Range("C4").Value = "招商惠润一年定期开放混合（MOM）A"
$q4.Range("D4").Value = "0.48"
$q4.Range("E4").Value = "68.20"
$q4.Range("F4").Value = "3.25"
$q4.Range("G4").Value = "0.0156"
$q4.Range("H4").Value = 4

# Row 5 - 001252 中海进取收益灵活配置混合
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "001252"
$q4.Range("C5").Value = "中海进取收益灵活配置混合"
$q4.Range("D5").Value = "0.51"
$q4.Range("E5").Value = "36.60"
$q4.Range("F5").Value = "1.53"
$q4.Range("G5").Value = "0.0078"
$q4.Range("H5").Value = 6

# Row 6 - 011215 招商惠润一年定期开放混合（MOM）C
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "011215"
$q4.Range("C6").Value = "招商惠润一年定期开放混合（MOM）C"
$q4.Range("D6").Value = "0.06"
$q4.Range("E6").Value = "68.20"
$q4.Range("F6").Value = "3.25"
$q4.Range("G6").Value = "0.0020"
$q4.Range("H6").Value = 4

# ---------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: insert a 2022-Q4 row before the
#    existing 2022-Q1 / 2021-Q4 rows, shifting those down by one row.
# ---------------------------------------------------------------------

# Make room: push current row 2 (2022-Q1) down to row 3, and add a brand
# new row 4 for 2021-Q4 (copy formats from row 3 downward first).
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.97

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.05
